$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D (Price) and E (Volume/1h) keep their original text
# representation (no float coercion / no scientific notation) by forcing
# a text number format before writing the literal strings back in.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "69.195.74"
$ws.Cells.Item(2, 5).Value = "  -3.60%  "
$ws.Cells.Item(3, 4).Value = "3.507.74"
$ws.Cells.Item(3, 5).Value = "  -5.00%  "
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).Value = "580.91"
$ws.Cells.Item(5, 5).Value = "  -1.51%  "
$ws.Cells.Item(6, 4).Value = "174.35"
$ws.Cells.Item(6, 5).Value = "  -3.56%  "
$ws.Cells.Item(7, 4).Value = "0.624"
$ws.Cells.Item(7, 5).Value = "  +0.63%  "
$ws.Cells.Item(8, 4).Value = "3.501.63"
$ws.Cells.Item(8, 5).Value = "  -4.92%  "
$ws.Cells.Item(9, 5).Value = "  +0.04%  "
$ws.Cells.Item(10, 5).Value = "  -6.12%  "
$ws.Cells.Item(11, 4).Value = "6.73"
$ws.Cells.Item(11, 5).Value = "  +6.50%  "
$ws.Cells.Item(12, 4).Value = "0.598"
$ws.Cells.Item(12, 5).Value = "  -2.54%  "
$ws.Cells.Item(13, 4).Value = "47.17"
$ws.Cells.Item(13, 5).Value = "  -6.03%  "
$ws.Cells.Item(14, 5).Value = "  -3.97%  "
$ws.Cells.Item(15, 4).Value = "677.30"
$ws.Cells.Item(15, 5).Value = "  -1.17%  "
$ws.Cells.Item(16, 4).Value = "4.070.10"
$ws.Cells.Item(17, 4).Value = "8.73"
$ws.Cells.Item(17, 5).Value = "  -3.49%  "
$ws.Cells.Item(18, 4).Value = "69.101.15"
$ws.Cells.Item(18, 5).Value = "  -3.85%  "
$ws.Cells.Item(19, 4).Value = "3.501.53"
$ws.Cells.Item(19, 5).Value = "  -5.20%  "
$ws.Cells.Item(20, 5).Value = "  -1.29%  "
$ws.Cells.Item(21, 4).Value = "17.50"
$ws.Cells.Item(21, 5).Value = "  -3.69%  "
$ws.Cells.Item(22, 4).Value = "11.20"
$ws.Cells.Item(22, 5).Value = "  -4.10%  "
$ws.Cells.Item(23, 4).Value = "0.904"
$ws.Cells.Item(23, 5).Value = "  -4.30%  "
$ws.Cells.Item(24, 4).Value = "16.17"
$ws.Cells.Item(24, 5).Value = "  -9.42%  "
$ws.Cells.Item(25, 4).Value = "98.05"
$ws.Cells.Item(25, 5).Value = "  -5.84%  "
$ws.Cells.Item(26, 4).Value = "3.87"
$ws.Cells.Item(26, 5).Value = "  -4.53%  "
$ws.Cells.Item(27, 5).Value = "  -0.22%  "
$ws.Cells.Item(28, 4).Value = "1.00"
$ws.Cells.Item(28, 5).Value = "  +0.04%  "
$ws.Cells.Item(29, 5).Value = "  -6.72%  "
$ws.Cells.Item(30, 4).Value = "9.45"
$ws.Cells.Item(30, 5).Value = "  -7.56%  "
$ws.Cells.Item(31, 4).Value = "32.95"
$ws.Cells.Item(31, 5).Value = "  -7.17%  "
$ws.Cells.Item(32, 4).Value = "8.73"
$ws.Cells.Item(32, 5).Value = "  -5.69%  "
$ws.Cells.Item(33, 4).Value = "3.21"
$ws.Cells.Item(33, 5).Value = "  -8.03%  "
$ws.Cells.Item(34, 5).Value = "  -6.18%  "
$ws.Cells.Item(35, 4).Value = "7.27"
$ws.Cells.Item(35, 5).Value = "  -1.22%  "
$ws.Cells.Item(36, 4).Value = "592.93"
$ws.Cells.Item(36, 5).Value = "  +4.43%  "
$ws.Cells.Item(37, 4).Value = "3.61"
$ws.Cells.Item(37, 5).Value = "  -15.44%  "
$ws.Cells.Item(38, 4).Value = "10.91"
$ws.Cells.Item(38, 5).Value = "  -3.88%  "
$ws.Cells.Item(39, 5).Value = "  -4.18%  "
$ws.Cells.Item(40, 5).Value = "  -3.84%  "
$ws.Cells.Item(41, 4).Value = "0.998"
$ws.Cells.Item(41, 5).Value = "  -0.04%  "
$ws.Cells.Item(42, 4).Value = "0.0440"
$ws.Cells.Item(42, 5).Value = "  -5.95%  "
$ws.Cells.Item(43, 4).Value = "0.337"
$ws.Cells.Item(43, 5).Value = "  -4.47%  "
$ws.Cells.Item(44, 2).Value = "Kaspa"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(44, 4).Value = "0.136"
$ws.Cells.Item(44, 5).Value = "  -7.90%  "
$ws.Cells.Item(45, 2).Value = "Maker"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(45, 4).Value = "3.421.60"
$ws.Cells.Item(45, 5).Value = "  -10.28%  "
$ws.Cells.Item(46, 4).Value = "33.44"
$ws.Cells.Item(46, 5).Value = "  -6.06%  "
$ws.Cells.Item(47, 4).Value = "0.0₃0708"
$ws.Cells.Item(47, 5).Value = "  -9.20%  "
$ws.Cells.Item(48, 4).Value = "2.92"
$ws.Cells.Item(48, 5).Value = "  -0.12%  "
$ws.Cells.Item(49, 4).Value = "2.60"
$ws.Cells.Item(49, 5).Value = "  -7.56%  "
$ws.Cells.Item(50, 5).Value = "  -0.58%  "
$ws.Cells.Item(51, 4).Value = "5.77"
$ws.Cells.Item(51, 5).Value = "  +17.96%  "
